$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.364.20'
$ws.Range('E2').Value = '  -5.06%  '
$ws.Range('D3').Value = '3.086.02'
$ws.Range('E3').Value = '  -5.25%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '''547.94'
$ws.Range('E5').Value = '  -6.19%  '
$ws.Range('D6').Value = '''136.37'
$ws.Range('E6').Value = '  -11.64%  '
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('D8').Value = '3.080.71'
$ws.Range('E8').Value = '  -5.15%  '
$ws.Range('D9').Value = '''0.494'
$ws.Range('E9').Value = '  -4.22%  '
$ws.Range('E10').Value = '  -6.15%  '
$ws.Range('D11').Value = '''6.27'
$ws.Range('E11').Value = '  -11.91%  '
$ws.Range('D12').Value = '''0.469'
$ws.Range('E12').Value = '  -4.50%  '
$ws.Range('D13').Value = '''35.27'
$ws.Range('E13').Value = '  -7.45%  '
$ws.Range('D14').Value = '''0.0000216'
$ws.Range('E14').Value = '  -8.51%  '
$ws.Range('D15').Value = '3.584.16'
$ws.Range('E15').Value = '  -5.28%  '
$ws.Range('D16').Value = '63.409.47'
$ws.Range('E16').Value = '  -5.05%  '
$ws.Range('D17').Value = '''0.111'
$ws.Range('E17').Value = '  -3.28%  '
$ws.Range('D18').Value = '3.092.22'
$ws.Range('E18').Value = '  -5.00%  '
$ws.Range('D19').Value = '''6.73'
$ws.Range('E19').Value = '  -5.95%  '
$ws.Range('D20').Value = '''486.24'
$ws.Range('E20').Value = '  -13.18%  '
$ws.Range('D21').Value = '''13.55'
$ws.Range('E21').Value = '  -6.78%  '
$ws.Range('D22').Value = '''0.715'
$ws.Range('E22').Value = '  -4.32%  '
$ws.Range('D23').Value = '''7.25'
$ws.Range('E23').Value = '  -7.07%  '
$ws.Range('D24').Value = '''78.90'
$ws.Range('E24').Value = '  -3.94%  '
$ws.Range('D25').Value = '''12.31'
$ws.Range('E25').Value = '  -10.01%  '
$ws.Range('D26').Value = '''0.999'
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('D27').Value = '''8.48'
$ws.Range('E27').Value = '  -8.96%  '
$ws.Range('D28').Value = '''2.75'
$ws.Range('E28').Value = '  -8.21%  '
$ws.Range('E29').Value = '  +0.11%  '
$ws.Range('D30').Value = '''1.96'
$ws.Range('E30').Value = '  -12.77%  '
$ws.Range('D31').Value = '''26.54'
$ws.Range('E31').Value = '  -5.12%  '
$ws.Range('E32').Value = '  -5.44%  '
$ws.Range('E33').Value = '  -10.13%  '
$ws.Range('D34').Value = '''59.11'
$ws.Range('E34').Value = '  +6.60%  '
$ws.Range('D35').Value = '''503.56'
$ws.Range('E35').Value = '  -10.75%  '
$ws.Range('D36').Value = '''6.04'
$ws.Range('E36').Value = '  -5.87%  '
$ws.Range('D37').Value = '''5.08'
$ws.Range('E37').Value = '  -11.42%  '
$ws.Range('D38').Value = '3.149.44'
$ws.Range('E38').Value = '  -0.92%  '
$ws.Range('D39').Value = '''0.0397'
$ws.Range('E39').Value = '  -13.76%  '
$ws.Range('D40').Value = '''0.0800'
$ws.Range('E40').Value = '  -7.67%  '
$ws.Range('E41').Value = '  -10.76%  '
$ws.Range('D42').Value = '''8.14'
$ws.Range('E42').Value = '  -6.08%  '
$ws.Range('D43').Value = '''2.61'
$ws.Range('E43').Value = '  -14.74%  '
$ws.Range('D44').Value = '''0.255'
$ws.Range('E44').Value = '  -7.82%  '
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('D46').Value = '''25.33'
$ws.Range('E46').Value = '  -4.56%  '
$ws.Range('E47').Value = '  -12.22%  '
$ws.Range('D48').Value = '''119.98'
$ws.Range('E48').Value = '  -5.03%  '
$ws.Range('E49').Value = '  -4.36%  '
$ws.Range('D50').Value = '0.0₃0505'
$ws.Range('E50').Value = '  -9.96%  '
$ws.Range('D51').Value = '''2.03'
$ws.Range('E51').Value = '  -9.96%  '
